# The presentation ships two themes:
#   ppt/theme/theme1.xml -> used by the slide master ("Simple Light" palette)
#   ppt/theme/theme2.xml -> used by the notes master ("Default" palette)
#
# The author's edit swaps the two themes' contents: the slide master's
# theme (theme1.xml) now carries the "Default" color palette (formerly
# theme2.xml's colors), while the notes master's theme (theme2.xml) would
# carry the former "Simple Light" palette. The font scheme / format scheme
# are identical between the two themes, so only the 12 theme colors change.
#
# Apply the new ("Default") color palette to the presentation's active
# theme (the slide master's theme) via the ThemeColorScheme COM API.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# Index -> (scheme slot, new RGB color as 0xBBGGRR COM long)
# 1  dk1       000000
# 2  lt1       FFFFFF
# 3  dk2       158158
# 4  lt2       F3F3F3
# 5  accent1   058DC7
# 6  accent2   50B432
# 7  accent3   ED561B
# 8  accent4   EDEF00
# 9  accent5   24CBE5
# 10 accent6   64E572
# 11 hlink     2200CC
# 12 folHlink  551A8B

$scheme.Colors(1).RGB  = 0
$scheme.Colors(2).RGB  = 16777215
$scheme.Colors(3).RGB  = 5800213
$scheme.Colors(4).RGB  = 15987699
$scheme.Colors(5).RGB  = 13077765
$scheme.Colors(6).RGB  = 3322960
$scheme.Colors(7).RGB  = 1791725
$scheme.Colors(8).RGB  = 61421
$scheme.Colors(9).RGB  = 15059748
$scheme.Colors(10).RGB = 7529828
$scheme.Colors(11).RGB = 13369378
$scheme.Colors(12).RGB = 9116245
